# Apply the changes described by the diff:
# 1. Update the "7.2.1" title strings to "7.2.1.1" in cells A1, B1, C1.
# 2. Update the selected cell in the sheet view from P9 to P7.
# 3. Set Q5 to 36.700000000000003 (was empty).
# 4. Update P6 from 13859.3 to 13859.2.
# 5. Update Q6 from 13979.1 to 13979.2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the title text in A1 / B1 / C1 (shared strings 12, 13, 14) ---
$ws.Range("A1").Value = " 7.2.1.1 Энергия керектөөлөрүнүн жалпы көлөмүндөгү энергиянын жаңыланган булактарынын  үлүшү"
$ws.Range("B1").Value = " 7.2.1.1 Доля возобновляемых источников энергии в общем объеме энергопотребления"
$ws.Range("C1").Value = "7.2.1.1 Renewable energy share in the total energy consumption"

# --- 2. Change the selected/active cell in the sheet view from P9 to P7 ---
$ws.Activate()
$ws.Range("P7").Select()

# --- 3. Fill in Q5 value ---
$ws.Range("Q5").Value = 36.700000000000003

# --- 4. / 5. Update P6 and Q6 values ---
$ws.Range("P6").Value = 13859.2
$ws.Range("Q6").Value = 13979.2
